$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore A2 to "a2" (it previously held the now-removed "Novo Valor" string).
# Assigning the same text as A1's already-shared "a2" string lets the saver
# reuse/compact the shared-strings table automatically.
$ws.Range("A2").Value = "a2"

# Widen column A (target stored width ~53.887 chars). ColumnWidth snaps to
# the sheet's pixel grid, so 53 reliably lands on the closest achievable
# stored width.
$ws.Columns.Item(1).ColumnWidth = 53
